$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.00", "5.32") are preserved exactly as authored, matching
# the source data which stores these as inline strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '56.985.16'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '2.325.85'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '523.97'
$ws.Range("E5").Value = '  +1.52%  '

$ws.Range("D6").Value = '132.41'
$ws.Range("E6").Value = '  -1.82%  '

$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.50%  '

$ws.Range("D8").Value = '0.534'
$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("D9").Value = '2.353.69'
$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("E10").Value = '  -1.31%  '

$ws.Range("E11").Value = '  +0.50%  '

$ws.Range("D12").Value = '5.32'
$ws.Range("E12").Value = '  -1.26%  '

$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  +0.47%  '

$ws.Range("D14").Value = '23.56'
$ws.Range("E14").Value = '  -1.51%  '

$ws.Range("D15").Value = '2.738.96'
$ws.Range("E15").Value = '  -0.31%  '

$ws.Range("D16").Value = '57.000.62'
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("D18").Value = '2.334.78'
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").Value = '336.56'
$ws.Range("E19").Value = '  +2.94%  '

$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("D21").Value = '6.95'
$ws.Range("E21").Value = '  +4.23%  '

$ws.Range("D22").Value = '4.16'
$ws.Range("E22").Value = '  -1.32%  '

$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").Value = '61.60'
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("D25").Value = '8.80'
$ws.Range("E25").Value = '  +10.46%  '

$ws.Range("D26").Value = '0.166'
$ws.Range("E26").Value = '  +0.46%  '

$ws.Range("D27").Value = '0.995'
$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("D28").Value = '1.34'
$ws.Range("E28").Value = '  +3.39%  '

$ws.Range("D29").Value = '169.36'
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").Value = '0.0₃0726'
$ws.Range("E31").Value = '  -1.78%  '

$ws.Range("D32").Value = '6.13'
$ws.Range("E32").Value = '  -1.37%  '

$ws.Range("D33").Value = '18.48'
$ws.Range("E33").Value = '  -0.38%  '

$ws.Range("D35").Value = '0.993'
$ws.Range("E35").Value = '  -0.42%  '

$ws.Range("D36").Value = '1.27'
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").Value = '0.914'
$ws.Range("E37").Value = '  +0.16%  '

$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '4.02'
$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("D39").Value = '1.59'
$ws.Range("E39").Value = '  +1.45%  '

$ws.Range("D40").Value = '38.91'
$ws.Range("E40").Value = '  +1.32%  '

$ws.Range("D41").Value = '148.98'
$ws.Range("E41").Value = '  +1.69%  '

$ws.Range("E42").Value = '  -1.89%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '3.60'
$ws.Range("E43").Value = '  -0.68%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '285.38'
$ws.Range("E44").Value = '  +2.25%  '

$ws.Range("D45").Value = '5.14'
$ws.Range("E45").Value = '  -0.39%  '

$ws.Range("D46").Value = '0.0932'
$ws.Range("E46").Value = '  -0.34%  '

$ws.Range("D47").Value = '0.0503'
$ws.Range("E47").Value = '  -0.37%  '

$ws.Range("D48").Value = '0.561'
$ws.Range("E48").Value = '  -0.23%  '

$ws.Range("D49").Value = '18.72'
$ws.Range("E49").Value = '  +4.04%  '

$ws.Range("D50").Value = '0.0217'
$ws.Range("E50").Value = '  -0.91%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '17.28'
$ws.Range("E51").Value = '  -1.74%  '
